# Fix input data reference errors:
# The "BCpUC" sheet computed a Battery Cost per Unit Capacity for year 2019,
# but the referenced lookup table (BBoSCpUC) only starts at year 2020, so the
# formula in B2 resolved to #N/A. Remove that bad 2019 row entirely so the
# remaining years shift up and the sheet recalculates cleanly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BCpUC")

$ws.Rows.Item(2).Delete()

$ws.Activate()
$ws.Range("E26").Select()
